# Final version edit pass:
#  - Slide 5 and Slide 8: mark as hidden ("Do not show"), and give Slide 8
#    the same slow, 2-second transition already used on Slide 5.
#  - Slide 24: merge the three consecutive runs of the "operators" bullet
#    into a single run (same text, same character formatting).

$p = $ppt.ActivePresentation

# --- Slide 5: hide it ---------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.SlideShowTransition.Hidden = 1

# --- Slide 8: hide it + add a slow / 2s transition -----------------------
$s8 = $p.Slides.Item(8)
$s8.SlideShowTransition.Hidden = 1
$s8.SlideShowTransition.Speed = 1
$s8.SlideShowTransition.Duration = 2

# --- Slide 24: merge the three runs of the 3rd bullet paragraph ----------
$s24 = $p.Slides.Item(24)
$shp = $s24.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3)

# paragraph currently reads: " " + "Modifications to the " + "operators " + "to move two examinations at the same time"
# keep the leading single-space run untouched; merge the remaining three
# runs (chars 2..73) into one run carrying the concatenated text.
$merged = $para.Characters(2, 72)
$merged.Text = "PLACEHOLDER_FOR_MERGE_TEMP"
$final = $para.Characters(2, 26)
$final.Text = "Modifications to the operators to move two examinations at the same time"
